$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new shift/date
$ws.Name = "Shift 2 - 2025-02-17"

# Update the activity log row: date, shift number, and engineer list.
# The Date column stores the date as literal text (not a real date value), so
# prefix with an apostrophe to keep Excel from reinterpreting the ISO-like
# string as a date serial when it's written back.
$ws.Range("A2").Value = "'2025-02-17"
$ws.Range("B2").Value = 2
$ws.Range("E2").Value = "Renato Hacel Cal y Mayor Rodríguez, Usuario de BC, Administrador"
